$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "birth_date"

$ws.Range("F2").Value = 29305
$ws.Range("F3").Value = 27943
$ws.Range("F4").Value = 31272
$ws.Range("F5").Value = 33134
$ws.Range("F6").Value = 33647
$ws.Range("F7").Value = 28516

$ws.Range("F2:F7").NumberFormat = "mm/dd/yy;@"
